$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 146.17053
$ws.Range("H2").Value = 438.51159
$ws.Range("I2").Value = 0.4047435297111188
$ws.Range("J2").Value = 0.4047435297111188
$ws.Range("M2").Value = 22.495411
$ws.Range("N2").Value = 67.486233
$ws.Range("O2").Value = 0.08292345339295874
$ws.Range("P2").Value = 0.08292345339295874
$ws.Range("Q2").Value = 3288.16614843783
$ws.Range("R2").Value = 29593.49533594047
$ws.Range("S2").Value = 0.03356273122210157
$ws.Range("T2").Value = 0.03356273122210157
$ws.Range("G3").Value = 146.17053
$ws.Range("H3").Value = 438.51159
$ws.Range("I3").Value = 0.4047435297111188
$ws.Range("J3").Value = 0.4047435297111188
$ws.Range("M3").Value = 82.64333833333332
$ws.Range("O3").Value = 0.3046430676248896
$ws.Range("P3").Value = 0.3046430676248896
$ws.Range("Q3").Value = 12080.02056515265
$ws.Range("R3").Value = 108720.1850863738
$ws.Range("S3").Value = 0.1233023104925209
$ws.Range("T3").Value = 0.1233023104925209
$ws.Range("G4").Value = 146.17053
$ws.Range("H4").Value = 438.51159
$ws.Range("I4").Value = 0.4047435297111188
$ws.Range("J4").Value = 0.4047435297111188
$ws.Range("M4").Value = 79.32606499999999
$ws.Range("N4").Value = 237.978195
$ws.Range("O4").Value = 0.2924148064631633
$ws.Range("P4").Value = 0.2924148064631633
$ws.Range("Q4").Value = 11595.13296386445
$ws.Range("R4").Value = 104356.19667478
$ws.Range("S4").Value = 0.1183530009076944
$ws.Range("T4").Value = 0.1183530009076944
$ws.Range("G5").Value = 146.17053
$ws.Range("H5").Value = 438.51159
$ws.Range("I5").Value = 0.4047435297111188
$ws.Range("J5").Value = 0.4047435297111188
$ws.Range("M5").Value = 14.467164
$ws.Range("N5").Value = 43.401492
$ws.Range("O5").Value = 0.05332941903938943
$ws.Range("P5").Value = 0.05332941903938943
$ws.Range("Q5").Value = 2114.67302947692
$ws.Range("R5").Value = 19032.05726529228
$ws.Range("S5").Value = 0.02158473729944582
$ws.Range("T5").Value = 0.02158473729944582
$ws.Range("G6").Value = 146.17053
$ws.Range("H6").Value = 438.51159
$ws.Range("I6").Value = 0.4047435297111188
$ws.Range("J6").Value = 0.4047435297111188
$ws.Range("M6").Value = 72.34725666666667
$ws.Range("N6").Value = 217.04177
$ws.Range("O6").Value = 0.2666892534795989
$ws.Range("P6").Value = 0.2666892534795989
$ws.Range("Q6").Value = 10575.0368510127
$ws.Range("R6").Value = 95175.3316591143
$ws.Range("S6").Value = 0.1079407497893561
$ws.Range("T6").Value = 0.1079407497893561
$ws.Range("H7").Value = 632.3552549999999
$ws.Range("I7").Value = 0.5836600531814327
$ws.Range("J7").Value = 0.5836600531814327
$ws.Range("M7").Value = 22.495411
$ws.Range("N7").Value = 67.486233
$ws.Range("O7").Value = 0.08292345339295874
$ws.Range("P7").Value = 0.08292345339295874
$ws.Range("Q7").Value = 4741.697119744935
$ws.Range("R7").Value = 42675.27407770441
$ws.Range("S7").Value = 0.04839910721732235
$ws.Range("T7").Value = 0.04839910721732235
$ws.Range("H8").Value = 632.3552549999999
$ws.Range("I8").Value = 0.5836600531814327
$ws.Range("J8").Value = 0.5836600531814327
$ws.Range("M8").Value = 82.64333833333332
$ws.Range("O8").Value = 0.3046430676248896
$ws.Range("P8").Value = 0.3046430676248896
$ws.Range("Q8").Value = 17419.98309527542
$ws.Range("S8").Value = 0.1778079890512979
$ws.Range("T8").Value = 0.1778079890512979
$ws.Range("H9").Value = 632.3552549999999
$ws.Range("I9").Value = 0.5836600531814327
$ws.Range("J9").Value = 0.5836600531814327
$ws.Range("M9").Value = 79.32606499999999
$ws.Range("N9").Value = 237.978195
$ws.Range("O9").Value = 0.2924148064631633
$ws.Range("P9").Value = 0.2924148064631633
$ws.Range("Q9").Value = 16720.75135374052
$ws.Range("R9").Value = 150486.7621836647
$ws.Range("S9").Value = 0.1706708414913282
$ws.Range("T9").Value = 0.1706708414913282
$ws.Range("H10").Value = 632.3552549999999
$ws.Range("I10").Value = 0.5836600531814327
$ws.Range("J10").Value = 0.5836600531814327
$ws.Range("M10").Value = 14.467164
$ws.Range("N10").Value = 43.401492
$ws.Range("O10").Value = 0.05332941903938943
$ws.Range("P10").Value = 0.05332941903938943
$ws.Range("Q10").Value = 3049.46239344894
$ws.Range("R10").Value = 27445.16154104046
$ws.Range("S10").Value = 0.03112625155266494
$ws.Range("T10").Value = 0.03112625155266494
$ws.Range("H11").Value = 632.3552549999999
$ws.Range("I11").Value = 0.5836600531814327
$ws.Range("J11").Value = 0.5836600531814327
$ws.Range("M11").Value = 72.34725666666667
$ws.Range("N11").Value = 217.04177
$ws.Range("O11").Value = 0.2666892534795989
$ws.Range("P11").Value = 0.2666892534795989
$ws.Range("Q11").Value = 15249.72264600015
$ws.Range("R11").Value = 137247.5038140013
$ws.Range("S11").Value = 0.1556558638688192
$ws.Range("T11").Value = 0.1556558638688193
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5494936666666667
$ws.Range("H12").Value = 1.648481
$ws.Range("I12").Value = 0.00152153793381314
$ws.Range("J12").Value = 0.00152153793381314
$ws.Range("M12").Value = 22.495411
$ws.Range("N12").Value = 67.486233
$ws.Range("O12").Value = 0.08292345339295874
$ws.Range("P12").Value = 0.08292345339295874
$ws.Range("Q12").Value = 12.36108587356367
$ws.Range("R12").Value = 111.249772862073
$ws.Range("S12").Value = 0.0001261711799401727
$ws.Range("T12").Value = 0.0001261711799401727
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5494936666666667
$ws.Range("H13").Value = 1.648481
$ws.Range("I13").Value = 0.00152153793381314
$ws.Range("J13").Value = 0.00152153793381314
$ws.Range("M13").Value = 82.64333833333332
$ws.Range("O13").Value = 0.3046430676248896
$ws.Range("P13").Value = 0.3046430676248896
$ws.Range("Q13").Value = 45.41199100635721
$ws.Range("R13").Value = 408.707919057215
$ws.Range("S13").Value = 0.0004635259836644713
$ws.Range("T13").Value = 0.0004635259836644713
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.5494936666666667
$ws.Range("H14").Value = 1.648481
$ws.Range("I14").Value = 0.00152153793381314
$ws.Range("J14").Value = 0.00152153793381314
$ws.Range("M14").Value = 79.32606499999999
$ws.Range("N14").Value = 237.978195
$ws.Range("O14").Value = 0.2924148064631633
$ws.Range("P14").Value = 0.2924148064631633
$ws.Range("Q14").Value = 43.58917031908832
$ws.Range("R14").Value = 392.302532871795
$ws.Range("S14").Value = 0.0004449202204423308
$ws.Range("T14").Value = 0.0004449202204423308
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.5494936666666667
$ws.Range("H15").Value = 1.648481
$ws.Range("I15").Value = 0.00152153793381314
$ws.Range("J15").Value = 0.00152153793381314
$ws.Range("M15").Value = 14.467164
$ws.Range("N15").Value = 43.401492
$ws.Range("O15").Value = 0.05332941903938943
$ws.Range("P15").Value = 0.05332941903938943
$ws.Range("Q15").Value = 7.949614992627999
$ws.Range("R15").Value = 71.546534933652
$ws.Range("S15").Value = 0.00008114273405664773
$ws.Range("T15").Value = 0.00008114273405664773
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.5494936666666667
$ws.Range("H16").Value = 1.648481
$ws.Range("I16").Value = 0.00152153793381314
$ws.Range("J16").Value = 0.00152153793381314
$ws.Range("M16").Value = 72.34725666666667
$ws.Range("N16").Value = 217.04177
$ws.Range("O16").Value = 0.2666892534795989
$ws.Range("P16").Value = 0.2666892534795989
$ws.Range("Q16").Value = 39.75435933904111
$ws.Range("R16").Value = 357.78923405137
$ws.Range("S16").Value = 0.0004057778157095177
$ws.Range("T16").Value = 0.0004057778157095177
$ws.Range("G17").Value = 3.410044
$ws.Range("H17").Value = 10.230132
$ws.Range("I17").Value = 0.009442349596941478
$ws.Range("J17").Value = 0.009442349596941478
$ws.Range("M17").Value = 22.495411
$ws.Range("N17").Value = 67.486233
$ws.Range("O17").Value = 0.08292345339295874
$ws.Range("P17").Value = 0.08292345339295874
$ws.Range("Q17").Value = 76.710341308084
$ws.Range("R17").Value = 690.3930717727559
$ws.Range("S17").Value = 0.0007829922367219995
$ws.Range("T17").Value = 0.0007829922367219995
$ws.Range("G18").Value = 3.410044
$ws.Range("H18").Value = 10.230132
$ws.Range("I18").Value = 0.009442349596941478
$ws.Range("J18").Value = 0.009442349596941478
$ws.Range("M18").Value = 82.64333833333332
$ws.Range("O18").Value = 0.3046430676248896
$ws.Range("P18").Value = 0.3046430676248896
$ws.Range("Q18").Value = 281.8174200235533
$ws.Range("R18").Value = 2536.35678021198
$ws.Range("S18").Value = 0.002876546346798892
$ws.Range("T18").Value = 0.002876546346798892
$ws.Range("G19").Value = 3.410044
$ws.Range("H19").Value = 10.230132
$ws.Range("I19").Value = 0.009442349596941478
$ws.Range("J19").Value = 0.009442349596941478
$ws.Range("M19").Value = 79.32606499999999
$ws.Range("N19").Value = 237.978195
$ws.Range("O19").Value = 0.2924148064631633
$ws.Range("P19").Value = 0.2924148064631633
$ws.Range("Q19").Value = 270.5053719968599
$ws.Range("R19").Value = 2434.54834797174
$ws.Range("S19").Value = 0.00276108282994717
$ws.Range("T19").Value = 0.00276108282994717
$ws.Range("G20").Value = 3.410044
$ws.Range("H20").Value = 10.230132
$ws.Range("I20").Value = 0.009442349596941478
$ws.Range("J20").Value = 0.009442349596941478
$ws.Range("M20").Value = 14.467164
$ws.Range("N20").Value = 43.401492
$ws.Range("O20").Value = 0.05332941903938943
$ws.Range("P20").Value = 0.05332941903938943
$ws.Range("Q20").Value = 49.33366579521599
$ws.Range("R20").Value = 444.0029921569439
$ws.Range("S20").Value = 0.000503555018371702
$ws.Range("T20").Value = 0.000503555018371702
$ws.Range("G21").Value = 3.410044
$ws.Range("H21").Value = 10.230132
$ws.Range("I21").Value = 0.009442349596941478
$ws.Range("J21").Value = 0.009442349596941478
$ws.Range("M21").Value = 72.34725666666667
$ws.Range("N21").Value = 217.04177
$ws.Range("O21").Value = 0.2666892534795989
$ws.Range("P21").Value = 0.2666892534795989
$ws.Range("Q21").Value = 246.7073285126266
$ws.Range("R21").Value = 2220.36595661364
$ws.Range("S21").Value = 0.002518173165101714
$ws.Range("T21").Value = 0.002518173165101715
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 0.228434
$ws.Range("H22").Value = 0.6853020000000001
$ws.Range("I22").Value = 0.0006325295766939459
$ws.Range("J22").Value = 0.0006325295766939459
$ws.Range("M22").Value = 22.495411
$ws.Range("N22").Value = 67.486233
$ws.Range("O22").Value = 0.08292345339295874
$ws.Range("P22").Value = 0.08292345339295874
$ws.Range("Q22").Value = 5.138716716374001
$ws.Range("R22").Value = 46.248450447366
$ws.Range("S22").Value = 0.00005245153687264834
$ws.Range("T22").Value = 0.00005245153687264834
$ws.Range("E23").Value = 3
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 0.228434
$ws.Range("H23").Value = 0.6853020000000001
$ws.Range("I23").Value = 0.0006325295766939459
$ws.Range("J23").Value = 0.0006325295766939459
$ws.Range("M23").Value = 82.64333833333332
$ws.Range("O23").Value = 0.3046430676248896
$ws.Range("P23").Value = 0.3046430676248896
$ws.Range("Q23").Value = 18.87854834883667
$ws.Range("R23").Value = 169.90693513953
$ws.Range("S23").Value = 0.0001926957506075166
$ws.Range("T23").Value = 0.0001926957506075166
$ws.Range("E24").Value = 3
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 0.228434
$ws.Range("H24").Value = 0.6853020000000001
$ws.Range("I24").Value = 0.0006325295766939459
$ws.Range("J24").Value = 0.0006325295766939459
$ws.Range("M24").Value = 79.32606499999999
$ws.Range("N24").Value = 237.978195
$ws.Range("O24").Value = 0.2924148064631633
$ws.Range("P24").Value = 0.2924148064631633
$ws.Range("Q24").Value = 18.12077033221
$ws.Range("R24").Value = 163.08693298989
$ws.Range("S24").Value = 0.0001849610137511868
$ws.Range("T24").Value = 0.0001849610137511868
$ws.Range("E25").Value = 3
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 0.228434
$ws.Range("H25").Value = 0.6853020000000001
$ws.Range("I25").Value = 0.0006325295766939459
$ws.Range("J25").Value = 0.0006325295766939459
$ws.Range("M25").Value = 14.467164
$ws.Range("N25").Value = 43.401492
$ws.Range("O25").Value = 0.05332941903938943
$ws.Range("P25").Value = 0.05332941903938943
$ws.Range("Q25").Value = 3.304792141176
$ws.Range("R25").Value = 29.743129270584
$ws.Range("S25").Value = 0.00003373243485031906
$ws.Range("T25").Value = 0.00003373243485031906
$ws.Range("E26").Value = 3
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 0.228434
$ws.Range("H26").Value = 0.6853020000000001
$ws.Range("I26").Value = 0.0006325295766939459
$ws.Range("J26").Value = 0.0006325295766939459
$ws.Range("M26").Value = 72.34725666666667
$ws.Range("N26").Value = 217.04177
$ws.Range("O26").Value = 0.2666892534795989
$ws.Range("P26").Value = 0.2666892534795989
$ws.Range("Q26").Value = 16.52657322939334
$ws.Range("R26").Value = 148.73915906454
$ws.Range("S26").Value = 0.0001686888406122751
$ws.Range("T26").Value = 0.0001686888406122752
